$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H16").Value = 6895
$ws.Range("J16").Value = 6895
$ws.Range("L16").Value = 6895
$ws.Range("N16").Value = -7355
$ws.Range("H32").Value = 857.9167
$ws.Range("J32").Value = 1042.1428
$ws.Range("L32").Value = 1042.1428
$ws.Range("N32").Value = -1694.1428
$ws.Range("H33").Value = 408.8889
$ws.Range("I33").Value = 311.42856
$ws.Range("J33").Value = 750
$ws.Range("K33").Value = 311.42856
$ws.Range("L33").Value = 750
$ws.Range("M33").Value = -82.42856
$ws.Range("N33").Value = -1208
$ws.Range("H43").Value = 1246
$ws.Range("I43").Value = 1053.6666
$ws.Range("K43").Value = 1053.6666
$ws.Range("M43").Value = -984.6666
$ws.Range("H92").Value = 794.8
$ws.Range("I92").Value = 849.75
$ws.Range("K92").Value = 849.75
$ws.Range("M92").Value = 398.25
$ws.Range("H98").Value = 1707.8182
$ws.Range("I98").Value = 1199
$ws.Range("K98").Value = 1199
$ws.Range("M98").Value = 299
$ws.Range("H100").Value = 2303.25
$ws.Range("J100").Value = 875
$ws.Range("L100").Value = 875
$ws.Range("N100").Value = -1957
$ws.Range("H101").Value = 308.33334
$ws.Range("I101").Value = 308.33334
$ws.Range("J101").Value = 0
$ws.Range("K101").Value = 925.0000200000001
$ws.Range("L101").Value = 0
$ws.Range("M101").Value = 696.9999799999999
$ws.Range("N101").ClearContents()
$ws.Range("H103").Value = 869
$ws.Range("J103").Value = 350
$ws.Range("L103").Value = 1050
$ws.Range("N103").Value = -2222
$ws.Range("H106").Value = 1549.75
$ws.Range("I106").Value = 1549.75
$ws.Range("K106").Value = 1549.75
$ws.Range("M106").Value = -918.75
$ws.Range("H107").Value = 1163.8948
$ws.Range("I107").Value = 860.16327
$ws.Range("K107").Value = 860.16327
$ws.Range("M107").Value = 1059.83673
$ws.Range("H113").Value = 7042
$ws.Range("I113").Value = 5557.7144
$ws.Range("K113").Value = 5557.7144
$ws.Range("M113").Value = -2303.7144
$ws.Range("H122").Value = 1707.8182
$ws.Range("I122").Value = 1199
$ws.Range("K122").Value = 3597
$ws.Range("M122").Value = -1147
$ws.Range("H132").Value = 4946.4614
$ws.Range("I132").Value = 4679.5654
$ws.Range("J132").Value = 6992.6665
$ws.Range("K132").Value = 14038.6962
$ws.Range("L132").Value = 20977.9995
$ws.Range("M132").Value = -11508.6962
$ws.Range("N132").Value = -26037.9995
$ws.Range("H137").Value = 2661.111
$ws.Range("I137").Value = 2158.3333
$ws.Range("K137").Value = 6474.999899999999
$ws.Range("M137").Value = -3924.999899999999
$ws.Range("H138").Value = 4656.5
$ws.Range("J138").Value = 4713.6553
$ws.Range("L138").Value = 14140.9659
$ws.Range("N138").Value = -24420.9659

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 20978.52
$ws.Range("I32").Value = 18519.291
$ws.Range("K32").Value = 18519.291
$ws.Range("M32").Value = -18232.291
$ws.Range("H97").Value = 1707.25
$ws.Range("I97").Value = 1707.25
$ws.Range("K97").Value = 1707.25
$ws.Range("M97").Value = -1211.25
$ws.Range("H113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("N113").ClearContents()
$ws.Range("H132").Value = 7009
$ws.Range("I132").Value = 2863.3333
$ws.Range("K132").Value = 8589.999899999999
$ws.Range("M132").Value = -6059.999899999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2118.818
$ws.Range("I86").Value = 2130.7
$ws.Range("K86").Value = 2130.7
$ws.Range("M86").Value = -1007.7
$ws.Range("H89").Value = 2118.818
$ws.Range("I89").Value = 2130.7
$ws.Range("K89").Value = 10653.5
$ws.Range("M89").Value = -5037.5
$ws.Range("H107").Value = 1829.2
$ws.Range("I107").Value = 1161.5
$ws.Range("K107").Value = 1161.5
$ws.Range("M107").Value = 758.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 5400
$ws.Range("I16").Value = 5400
$ws.Range("K16").Value = 5400
$ws.Range("M16").Value = -5113
$ws.Range("H25").Value = 552.7778
$ws.Range("I25").Value = 552.7778
$ws.Range("K25").Value = 552.7778
$ws.Range("M25").Value = -378.7778
$ws.Range("H51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("N51").ClearContents()
$ws.Range("H61").Value = 0
$ws.Range("J61").Value = 0
$ws.Range("L61").Value = 0
$ws.Range("N61").ClearContents()
$ws.Range("H80").Value = 0
$ws.Range("I80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("M80").ClearContents()
$ws.Range("H83").Value = 0
$ws.Range("I83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("M83").ClearContents()
$ws.Range("H103").Value = 16674.334
$ws.Range("J103").Value = 0
$ws.Range("L103").Value = 0
$ws.Range("N103").ClearContents()
$ws.Range("H107").Value = 560.3125
$ws.Range("I107").Value = 453.7
$ws.Range("K107").Value = 453.7
$ws.Range("M107").Value = 1466.3
$ws.Range("H113").Value = 5400
$ws.Range("I113").Value = 5400
$ws.Range("K113").Value = 5400
$ws.Range("M113").Value = -3230
$ws.Range("H122").Value = 3041.7222
$ws.Range("I122").Value = 2926.5293
$ws.Range("K122").Value = 8779.5879
$ws.Range("M122").Value = -6329.5879

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 853.7273
$ws.Range("I23").Value = 750
$ws.Range("J23").Value = 913
$ws.Range("K23").Value = 2250
$ws.Range("L23").Value = 2739
$ws.Range("M23").Value = -2015
$ws.Range("N23").Value = -3209
$ws.Range("H86").Value = 666.6667
$ws.Range("I86").Value = 300
$ws.Range("K86").Value = 900
$ws.Range("M86").Value = 286
$ws.Range("H89").Value = 666.6667
$ws.Range("I89").Value = 300
$ws.Range("K89").Value = 2700
$ws.Range("M89").Value = 3228
$ws.Range("H97").Value = 1320.1538
$ws.Range("I97").Value = 1576
$ws.Range("J97").Value = 1160.25
$ws.Range("K97").Value = 4728
$ws.Range("L97").Value = 3480.75
$ws.Range("M97").Value = -4232
$ws.Range("N97").Value = -4472.75
$ws.Range("H137").Value = 22606
$ws.Range("I137").Value = 20515
$ws.Range("K137").Value = 61545
$ws.Range("M137").Value = -56445

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2013.1428
$ws.Range("I7").Value = 2013.1428
$ws.Range("K7").Value = 2013.1428
$ws.Range("M7").Value = -1901.1428
$ws.Range("H10").Value = 493.5
$ws.Range("J10").Value = 493.5
$ws.Range("L10").Value = 493.5
$ws.Range("N10").Value = -773.5
$ws.Range("H22").Value = 758.4666999999999
$ws.Range("I22").Value = 762.6429000000001
$ws.Range("K22").Value = 762.6429000000001
$ws.Range("M22").Value = -467.6429000000001
$ws.Range("H27").Value = 758.4666999999999
$ws.Range("I27").Value = 762.6429000000001
$ws.Range("K27").Value = 762.6429000000001
$ws.Range("M27").Value = -655.6429000000001
$ws.Range("H46").Value = 2239
$ws.Range("I46").Value = 2239
$ws.Range("K46").Value = 2239
$ws.Range("M46").Value = -2051
$ws.Range("H47").Value = 0
$ws.Range("J47").Value = 0
$ws.Range("L47").Value = 0
$ws.Range("N47").ClearContents()
$ws.Range("H52").Value = 0
$ws.Range("J52").Value = 0
$ws.Range("L52").Value = 0
$ws.Range("N52").ClearContents()
$ws.Range("H61").Value = 3764.7
$ws.Range("I61").Value = 3516.3333
$ws.Range("K61").Value = 3516.3333
$ws.Range("M61").Value = -3314.3333
$ws.Range("H113").Value = 3764.7
$ws.Range("I113").Value = 3516.3333
$ws.Range("K113").Value = 3516.3333
$ws.Range("M113").Value = -1346.3333
$ws.Range("H126").Value = 2013.1428
$ws.Range("I126").Value = 2013.1428
$ws.Range("K126").Value = 6039.428400000001
$ws.Range("M126").Value = -3569.428400000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 340.81818
$ws.Range("I107").Value = 331.25
$ws.Range("K107").Value = 993.75
$ws.Range("M107").Value = 926.25
$ws.Range("H132").Value = 3484.5
$ws.Range("I132").Value = 3381.4
$ws.Range("K132").Value = 10144.2
$ws.Range("M132").Value = -7614.200000000001
